$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current "Dtem" column (E), shifting
# Dtem..notes from E:L to F:M.
$ws.Columns("E:E").Insert()

# New column keeps the same width as the Ligand column (D) next to it.
$ws.Range("E1").ColumnWidth = $ws.Range("D1").ColumnWidth

# Ligand column (D) data update: Carboxylate -> Amine for the Polystyrene rows.
$ws.Range("D3").Value = "Amine"
$ws.Range("D4").Value = "Amine"
$ws.Range("D5").Value = "Amine"
$ws.Range("D6").Value = "Amine"

# Header for the newly inserted column.
$ws.Range("E1").Value = "BET"

# Zeta Potential column (B) data updates.
$ws.Range("B3").Value = 100
$ws.Range("B4").Formula = "=AVERAGE(-19,-26,-28)"
$ws.Range("B5").Formula = "=AVERAGE(-19,-26,-28)"
$ws.Range("B6").Value = 9

# Restore the selection to match the saved workbook state.
[void]$ws.Range("E2").Select()
